$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# ---------------------------------------------------------------
# Helper: append a new run of text right before a paragraph's own
# end-of-paragraph mark (so it becomes the last run in that
# paragraph). Returns the Range of the newly inserted text so the
# caller can tweak Font properties (color, italics, ...).
# ---------------------------------------------------------------
function Append-RunToCell($row, $col, $text) {
    $cell = $t.Cell($row, $col)
    $p = $cell.Range.Paragraphs.Item(1)
    $insertPos = $p.Range.End - 1
    $insertPoint = $d.Range($insertPos, $insertPos)
    $insertPoint.InsertAfter($text)
    $newRange = $d.Range($insertPos, $insertPos + $text.Length)
    return $newRange
}

# ---------------------------------------------------------------
# Helper: replace the entire text of a (single paragraph) cell with
# new text, keeping the formatting of the first original run.
# ---------------------------------------------------------------
function Set-CellText($row, $col, $text) {
    $cell = $t.Cell($row, $col)
    $p = $cell.Range.Paragraphs.Item(1)
    $rng = $d.Range($p.Range.Start, $p.Range.End - 1)
    $rng.Text = $text
}

# =================================================================
# 1. "Writing functions" -> append trailing space
# =================================================================
Append-RunToCell 3 3 " " | Out-Null

# =================================================================
# 2. "General structure of a simulation" -> append " Part 1"
# =================================================================
Append-RunToCell 4 3 " Part 1" | Out-Null

# =================================================================
# 3. "Factorial vs. one-at-at-time simulations" ->
#    "General structure of a simulation Part 2 [assignment]"
#    (last piece, "[assignment]", is colored red)
# =================================================================
Set-CellText 5 3 "General structure of a simulation Part 2 "
$assignRange3 = Append-RunToCell 5 3 "[assignment]"
$assignRange3.Font.Color = 255

# =================================================================
# 4. "Understanding p-values" (3 runs) ->
#    "The impact of violating statistical assumptions" (1 run)
# =================================================================
Set-CellText 6 3 "The impact of violating statistical assumptions"

# =================================================================
# 5. "Hidden multiplicity in ANOVA" -> append " [assignment]" (red)
# =================================================================
Append-RunToCell 7 3 " " | Out-Null
$assignRange7 = Append-RunToCell 7 3 "[assignment]"
$assignRange7.Font.Color = 255

# =================================================================
# 6. "What does it mean to violate assumptions?" ->
#    "Is it worth testing statistical assumptions?"
# =================================================================
Set-CellText 8 3 "Is it worth testing statistical assumptions?"

# =================================================================
# 7. Row 9: delete the 2nd paragraph ("Otherwise: Simulating causal
#    models"), keep the first ("<<Probably no class ...>>")
# =================================================================
$cell9 = $t.Cell(9, 3)
$cell9.Range.Paragraphs.Item(2).Range.Delete()

# =================================================================
# 8/9. Rows 10 & 12 swap their text content. Row 10 additionally
#      gains szCs=22 on the paragraph-mark run props, and a trailing
#      space run.
# =================================================================
Set-CellText 10 3 "Understanding Confidence Intervals via sequential testing"
Append-RunToCell 10 3 " " | Out-Null
$cell10 = $t.Cell(10, 3)
$p10 = $cell10.Range.Paragraphs.Item(1)
$p10.Range.Font.SizeBi = 11

Set-CellText 12 3 "The difference between significant and non-significant is not itself significant"

# =================================================================
# 11/12/13. Rows 13, 14, 15 shift down (13<-14, 14<-15, 15<-new).
#           Row 13 additionally gains szCs=22 on its run.
# =================================================================
Set-CellText 13 3 "How standardized are ‘standardized’ effect sizes?"
$cell13 = $t.Cell(13, 3)
$p13 = $cell13.Range.Paragraphs.Item(1)
$p13.Range.Font.SizeBi = 11

Set-CellText 14 3 "Meta-analysis and bias"

Set-CellText 15 3 "Simulating causal models"

Write-Output "done"
